# "Generate Report for Archive"
# Refresh the localization status report:
#   - the two handback/translation files have moved from "Ready for handoff"
#     to "In Translation"
#   - the (now shorter) status text lets the status columns shrink a bit

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the status columns to fit the new, shorter text ---
# (was ~17.22 characters wide, now ~13.41 characters wide)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5       # column C (Status)
